# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp
# - Reorder two pairs/blocks of countries in the ranking (their case-count
#   rows swapped position in the source feed) and refresh the underlying
#   case statistics (columns B:H) for every affected row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 13 de Agosto de 2020 a las 17:34"

# --- Estados Unidos (row 4) --------------------------------------------
$ws.Range("B4").Value = 5365495
$ws.Range("C4").Value = 5193
$ws.Range("D4").Value = 2813845
$ws.Range("E4").Value = 2382425
$ws.Range("G4").Value = 94
$ws.Range("H4").Value = 169225

# --- India (row 6) ------------------------------------------------------
$ws.Range("B6").Value = 2431558
$ws.Range("C6").Value = 36087
$ws.Range("D6").Value = 1725834
$ws.Range("E6").Value = 658197
$ws.Range("G6").Value = 389
$ws.Range("H6").Value = 47527

# --- Italia (row 20) ------------------------------------------------------
$ws.Range("B20").Value = 252235
$ws.Range("C20").Value = 522
$ws.Range("D20").Value = 202923
$ws.Range("E20").Value = 14081
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 35231

# --- Oman / Republica Dominicana swap (rows 37-38) -----------------------
$ws.Range("A37").Value = "Republica Dominicana"
$ws.Range("B37").Value = 83134
$ws.Range("C37").Value = 910
$ws.Range("D37").Value = 47946
$ws.Range("E37").Value = 33795
$ws.Range("G37").Value = 22
$ws.Range("H37").Value = 1393

$ws.Range("A38").Value = "Oman"
$ws.Range("B38").Value = 82531
$ws.Range("C38").Value = 232
$ws.Range("D38").Value = 77278
$ws.Range("E38").Value = 4702
$ws.Range("G38").Value = 12
$ws.Range("H38").Value = 551

# --- Bielorrusia (row 42) -------------------------------------------------
$ws.Range("B42").Value = 69203
$ws.Range("C42").Value = 101
$ws.Range("D42").Value = 66178
$ws.Range("E42").Value = 2426
$ws.Range("G42").Value = 4
$ws.Range("H42").Value = 599

# --- Guatemala (row 49) -------------------------------------------------
$ws.Range("B49").Value = 53548
$ws.Range("C49").Value = 325
$ws.Range("E49").Value = 12838
$ws.Range("G49").Value = 6
$ws.Range("H49").Value = 1770

# --- row 53 ----------------------------------------------------------------
$ws.Range("E53").Value = 3261
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 167

# --- row 61 ------------------------------------------------------------
$ws.Range("B61").Value = 33915
$ws.Range("C61").Value = 91
$ws.Range("D61").Value = 31269
$ws.Range("E61").Value = 2146
$ws.Range("G61").Value = 3
$ws.Range("H61").Value = 500

# --- row 64 ------------------------------------------------------------
$ws.Range("D64").Value = 20276
$ws.Range("E64").Value = 7933

# --- row 99 ------------------------------------------------------------
$ws.Range("B99").Value = 6971
$ws.Range("C99").Value = 154
$ws.Range("D99").Value = 3616
$ws.Range("E99").Value = 3142
$ws.Range("G99").Value = 5
$ws.Range("H99").Value = 213

# --- row 120 -----------------------------------------------------------
$ws.Range("B120").Value = 3119
$ws.Range("C120").Value = 28
$ws.Range("D120").Value = 2940
$ws.Range("E120").Value = 140

# --- Trinidad y Tobago / Birmania / Guadalupe / Islas Feroe / Mauricio /
#     Martinica / Isla de Man rotation (rows 170-176) ---------------------
$ws.Range("A170").Value = "Trinidad yTobago"
$ws.Range("B170").Value = 369
$ws.Range("C170").Value = 43
$ws.Range("D170").Value = 139
$ws.Range("E170").Value = 222
$ws.Range("H170").Value = 8

$ws.Range("A171").Value = "Birmania"
$ws.Range("B171").Value = 369
$ws.Range("C171").Value = 8
$ws.Range("D171").Value = 321
$ws.Range("E171").Value = 42
$ws.Range("H171").Value = 6

$ws.Range("A172").Value = "Guadalupe"
$ws.Range("B172").Value = 367
$ws.Range("D172").Value = 289
$ws.Range("E172").Value = 64
$ws.Range("H172").Value = 14

$ws.Range("A173").Value = "Islas Feroe"
$ws.Range("B173").Value = 362
$ws.Range("C173").Value = 23
$ws.Range("D173").Value = 225
$ws.Range("E173").Value = 137
$ws.Range("H173").Value = 0

$ws.Range("A174").Value = "Mauricio"
$ws.Range("B174").Value = 344
$ws.Range("D174").Value = 334
$ws.Range("E174").Value = 0
$ws.Range("H174").Value = 10

$ws.Range("A175").Value = "Martinica"
$ws.Range("D175").Value = 98
$ws.Range("E175").Value = 222
$ws.Range("H175").Value = 16

$ws.Range("A176").Value = "Isla de Man"
$ws.Range("B176").Value = 336
$ws.Range("D176").Value = 312
$ws.Range("E176").Value = 0
$ws.Range("H176").Value = 24
